$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Fix the typo "realtive" -> "relative" in the opening paragraph.
#    Match only up through "real" (i.e. stop right at the existing
#    "_GoBack" bookmark) so the bookmark itself is not swallowed by
#    the replace.
# ------------------------------------------------------------------
$d.Content.Find.Execute("We highlight the most prevalent trends real", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "We highlight the most prevalent trends rela", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Correct the inflated respondent count "1100000" -> "1100".
# ------------------------------------------------------------------
$d.Content.Find.Execute("This report includes insights from 1100000", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "This report includes insights from 1100", 2) | Out-Null

# ------------------------------------------------------------------
# 3) Correct the inflated industry count "3500" -> "35".
# ------------------------------------------------------------------
$d.Content.Find.Execute("representing more than 3500", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "representing more than 35", 2) | Out-Null

# ------------------------------------------------------------------
# 4) Move the "_GoBack" bookmark from its old position (between
#    "real" and "tive" in paragraph 1) to sit right after
#    "more than 35" in paragraph 2 - matching where it ends up after
#    the edit.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$full = $d.Content.Text
$idx = $full.IndexOf("more than 35")
$pos = $idx + ("more than 35").Length
$d.Bookmarks.Add("_GoBack", $d.Range($pos, $pos)) | Out-Null

Write-Output $d.Paragraphs.Item(2).Range.Text
Write-Output $d.Paragraphs.Item(3).Range.Text
